$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Input): remove the "office" row ---
$ws1.Rows.Item(1).Delete()

# After the delete, rows are now:
# 1 dateofbirth(date-style) / 2 firstname / 3 middlename / 4 lastname / 5 active / 6 activationdate / 7 submittedon
# Target order: firstname, middlename, lastname, dateofbirth, active, activationdate, submittedon
# Swap the B-column formatting of row1 (date style) and row4 (text style) so it follows the content move.
# Use untouched reference cells of each style (row3 = text style, row7 = date style) as copy sources.
$ws1.Cells.Item(3,2).Copy()
$ws1.Cells.Item(1,2).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(7,2).Copy()
$ws1.Cells.Item(4,2).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Capture current values for the 4 rows being reordered
$dobA = $ws1.Cells.Item(1,1).Value2
$dobB = $ws1.Cells.Item(1,2).Value2
$fnA  = $ws1.Cells.Item(2,1).Value2
$mnA  = $ws1.Cells.Item(3,1).Value2
$mnB  = $ws1.Cells.Item(3,2).Value2
$lnA  = $ws1.Cells.Item(4,1).Value2

# Write the new order with updated first/last name values
$ws1.Cells.Item(1,1).Value = $fnA
$ws1.Cells.Item(1,2).Value = "Jhon"
$ws1.Cells.Item(2,1).Value = $mnA
$ws1.Cells.Item(2,2).Value = $mnB
$ws1.Cells.Item(3,1).Value = $lnA
$ws1.Cells.Item(3,2).Value = "Deer"
$ws1.Cells.Item(4,1).Value = $dobA
$ws1.Cells.Item(4,2).Value = $dobB

# Sheet1 view: new selection
$ws1.Range("E14").Select() | Out-Null

# --- Sheet2 (Output): update the computed display name ---
$ws2.Cells.Item(1,2).Value = "Jhon N Deer"
$ws2.Range("C14").Select() | Out-Null

# Activate Input sheet last so it is the active tab in the saved workbook
$ws1.Activate() | Out-Null
